# Update InsideBet Data: Automatizado
# The fixture that previously occupied row 30 (Wk24, Fri 2026-02-20,
# Fortuna Sittard vs Excelsior) has dropped off the upcoming-fixtures list.
# Remove that row entirely so every following fixture shifts up by one row
# (matching the new dimension A1:L137).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(30).Delete()
